$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.318.27"
$ws.Range("E2").Value = '  +1.14%  '
$ws.Range("D3").Value = "'1.913.28"
$ws.Range("E3").Value = '  +1.50%  '
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = "'320.78"
$ws.Range("E5").Value = '  -2.94%  '
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("D7").Value = "'0.4729"
$ws.Range("E7").Value = '  +3.17%  '
$ws.Range("D8").Value = "'0.4069"
$ws.Range("E8").Value = '  +0.41%  '
$ws.Range("D9").Value = "'47.76"
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("D10").Value = "'0.08041"
$ws.Range("E10").Value = '  +1.01%  '
$ws.Range("D11").Value = "'1.002"
$ws.Range("E11").Value = '  +1.19%  '
$ws.Range("D12").Value = "'22.77"
$ws.Range("E12").Value = '  +5.32%  '
$ws.Range("D13").Value = "'1.913.59"
$ws.Range("E13").Value = '  +1.60%  '
$ws.Range("D14").Value = "'5.892"
$ws.Range("E14").Value = '  -0.15%  '
$ws.Range("D15").Value = "'7.111"
$ws.Range("E15").Value = '  +0.65%  '
$ws.Range("D16").Value = "'89.59"
$ws.Range("E16").Value = '  +1.54%  '
$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = '  -0.05%  '
$ws.Range("D18").Value = "'0.06622"
$ws.Range("E18").Value = '  +1.08%  '
$ws.Range("D19").Value = "'0.00001030"
$ws.Range("E19").Value = '  +0.10%  '
$ws.Range("D20").Value = "'17.68"
$ws.Range("E20").Value = '  +1.50%  '
$ws.Range("D21").Value = "'0.9999"
$ws.Range("E21").Value = '  +0.01%  '
$ws.Range("D22").Value = "'29.336.65"
$ws.Range("E22").Value = '  +1.14%  '
$ws.Range("D23").Value = "'5.519"
$ws.Range("E23").Value = '  +2.04%  '
$ws.Range("D24").Value = "'11.47"
$ws.Range("E24").Value = '  +0.85%  '
$ws.Range("D25").Value = "'2.202"
$ws.Range("E25").Value = '  -0.18%  '
$ws.Range("D26").Value = "'2.139.16"
$ws.Range("E26").Value = '  +2.55%  '
$ws.Range("D27").Value = "'153.84"
$ws.Range("E27").Value = '  -1.73%  '
$ws.Range("D28").Value = "'19.78"
$ws.Range("E28").Value = '  +1.03%  '
$ws.Range("D29").Value = "'6.029"
$ws.Range("E29").Value = '  +11.35%  '
$ws.Range("E30").Value = '  +0.57%  '
$ws.Range("D31").Value = "'117.67"
$ws.Range("E31").Value = '  -0.07%  '
$ws.Range("D32").Value = "'1.077"
$ws.Range("E32").Value = '  +7.26%  '
$ws.Range("D33").Value = "'0.09510"
$ws.Range("E33").Value = '  +1.97%  '
$ws.Range("E34").Value = '  +1.38%  '
$ws.Range("D35").Value = "'3.545"
$ws.Range("E35").Value = '  -1.46%  '
$ws.Range("D36").Value = "'5.388"
$ws.Range("E36").Value = '  +2.17%  '
$ws.Range("D37").Value = "'0.06081"
$ws.Range("E37").Value = '  +0.56%  '
$ws.Range("D38").Value = "'0.02249"
$ws.Range("E38").Value = '  +1.57%  '
$ws.Range("D39").Value = "'8.252"
$ws.Range("E39").Value = '  +0.04%  '
$ws.Range("D40").Value = "'1.174"
$ws.Range("E40").Value = '  +0.13%  '
$ws.Range("D41").Value = "'0.5858"
$ws.Range("E41").Value = '  +1.54%  '
$ws.Range("D42").Value = "'2.515"
$ws.Range("E42").Value = '  +11.38%  '
$ws.Range("E43").Value = '  +0.84%  '
$ws.Range("D44").Value = "'10.13"
$ws.Range("E44").Value = '  +0.32%  '
$ws.Range("D45").Value = "'0.07904"
$ws.Range("E45").Value = '  +5.09%  '
$ws.Range("E46").Value = '  +1.29%  '
$ws.Range("D47").Value = "'12.17"
$ws.Range("E47").Value = '  +1.04%  '
$ws.Range("D48").Value = "'0.5510"
$ws.Range("E48").Value = '  +1.22%  '
$ws.Range("D49").Value = "'1.923"
$ws.Range("E49").Value = '  +1.49%  '
$ws.Range("E50").Value = '  +2.05%  '
$ws.Range("D51").Value = "'44.25"
$ws.Range("E51").Value = '  -2.12%  '

# Reset number formatting/style on text-forced D-column cells to match original (no explicit style)
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
